$wb = $excel.ActiveWorkbook

# --- Bus sheet: remove "bus 1" row (row 3) ---
$wsBus = $wb.Worksheets.Item("Bus")
$wsBus.Rows.Item(3).Delete()
$wsBus.Range("B7").Select()

# --- Load sheet: remove "demand 1" row (row 3) ---
$wsLoad = $wb.Worksheets.Item("Load")
$wsLoad.Rows.Item(3).Delete()
$wsLoad.Range("A3:XFD3").Select()

# --- Line sheet: remove "line_0-1" row (row 2) ---
$wsLine = $wb.Worksheets.Item("Line")
$wsLine.Rows.Item(2).Delete()
$wsLine.Range("B8").Select()

# --- Generator sheet: wind bus 1 -> bus 0, and False -> True for wind & diesel ---
$wsGen = $wb.Worksheets.Item("Generator")
$wsGen.Range("C2").Value = "bus 0"
$wsGen.Range("D2").Value = "True"
$wsGen.Range("D3").Value = "True"
$wsGen.Range("D4").Select()

# --- StorageUnit sheet: remove "hydro" row (row 2) ---
$wsStorage = $wb.Worksheets.Item("StorageUnit")
$wsStorage.Rows.Item(2).Delete()
$wsStorage.Range("B9").Select()
